$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.667.86'
$ws.Range('E2').Value = '  -6.67%  '
$ws.Range('D3').Value = '2.446.36'
$ws.Range('E3').Value = '  -9.80%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '466.73'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -7.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.37'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.75%  '
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.495'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -6.35%  '
$ws.Range('D9').Value = '2.444.50'
$ws.Range('E9').Value = '  -10.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0958'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -9.00%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.32'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -12.01%  '
$ws.Range('E12').Value = '  -9.30%  '
$ws.Range('E13').Value = '  -3.73%  '
$ws.Range('D14').Value = '2.874.60'
$ws.Range('E14').Value = '  -9.81%  '
$ws.Range('D15').Value = '54.721.86'
$ws.Range('E15').Value = '  -6.65%  '
$ws.Range('E16').Value = '  -0.47%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.75'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -8.59%  '
$ws.Range('D18').Value = '2.448.59'
$ws.Range('E18').Value = '  -10.43%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.23'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -11.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '312.84'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -8.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.58'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -12.47%  '
$ws.Range('E22').Value = '  +0.34%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.69'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.88%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.38'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -13.89%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '56.71'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -10.12%  '
$ws.Range('E26').Value = '  +1.23%  '
$ws.Range('B27').Value = 'Polygon'
$ws.Range('C27').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.387'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -9.17%  '
$ws.Range('B28').Value = 'Kaspa'
$ws.Range('C28').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.157'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -9.49%  '
$ws.Range('D29').Value = '2.542.75'
$ws.Range('E29').Value = '  -10.48%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.17'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.15%  '
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('D32').Value = '0.0₃0721'
$ws.Range('E32').Value = '  -12.46%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '147.01'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.88%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.80'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -7.28%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.43'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -10.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.03'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -7.17%  '
$ws.Range('E37').Value = '  -14.24%  '
$ws.Range('E38').Value = '  -6.13%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.802'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -14.63%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '32.99'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -8.58%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.597'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('E43').Value = '  -5.80%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.26'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -8.81%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.24'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -10.30%  '
$ws.Range('E46').Value = '  -2.67%  '
$ws.Range('D47').Value = '1.943.12'
$ws.Range('E47').Value = '  -11.29%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0884'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.52%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0218'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.19%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '234.47'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.17%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '16.64'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -12.13%  '
